$p = $ppt.ActivePresentation

# New text for each named rectangle inside the "Group 3" diagram on the
# Android architecture slide.
$replacements = @{
    "Rectangle 7"  = "Linux Application"
    "Rectangle 8"  = "/dev"
    "Rectangle 9"  = "/proc"
    "Rectangle 10" = "dmesg"
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $top = $s.Shapes.Item($shi)
        if ($top.Name -ne "Group 3") { continue }

        for ($i = 1; $i -le $top.GroupItems.Count; $i++) {
            $shp = $top.GroupItems.Item($i)
            if (-not $replacements.ContainsKey($shp.Name)) { continue }
            if (-not $shp.HasTextFrame) { continue }
            if (-not $shp.TextFrame.HasText) { continue }

            $shp.TextFrame.TextRange.Text = $replacements[$shp.Name]
        }
    }
}
